# Restore selection on Feuil1 (existing sheet) before switching away
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B13").Select()

# Add the new "Feuil2" worksheet after the existing "Feuil1"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "Feuil2"

# Populate the quizz questions / answers / results
$ws2.Cells.Item(1, 1).Value = 'INTITULE QUESTION'
$ws2.Cells.Item(1, 2).Value = 'REPONSE'
$ws2.Cells.Item(1, 3).Value = 'RESULTAT'
$ws2.Cells.Item(2, 1).Value = 'Le Parlement français est composé de deux chambres, l''Assemblée nationale et le Sénat qui ont respectivement, jusqu''en 2010'
$ws2.Cells.Item(2, 2).Value = '577 députés et 343 sénateurs'
$ws2.Cells.Item(2, 3).Value = 1
$ws2.Cells.Item(3, 2).Value = '571 députés et 338 sénateurs'
$ws2.Cells.Item(3, 3).Value = 0
$ws2.Cells.Item(4, 2).Value = '574 députés et 335 sénateurs'
$ws2.Cells.Item(4, 3).Value = 0
$ws2.Cells.Item(5, 2).Value = '581 députés et 348 sénateurs'
$ws2.Cells.Item(5, 3).Value = 0
$ws2.Cells.Item(6, 1).Value = 'Les deux chambres du Parlement siègent respectivement'
$ws2.Cells.Item(6, 2).Value = 'Au Palais Bourbon pour les députés, au Palais du Luxembourg pour les sénateurs'
$ws2.Cells.Item(6, 3).Value = 1
$ws2.Cells.Item(7, 2).Value = 'Au Palais du Luxembourg pour les députés, au Palais Bourbon pour les sénateurs'
$ws2.Cells.Item(7, 3).Value = 0
$ws2.Cells.Item(8, 2).Value = 'Au Palais Royal pour les députés, au Palais du Luxembourg pour les sénateurs'
$ws2.Cells.Item(8, 3).Value = 0
$ws2.Cells.Item(9, 2).Value = 'Au Grand Palais pour les députés, au Palais Royal pour les sénateurs'
$ws2.Cells.Item(9, 3).Value = 0
$ws2.Cells.Item(10, 1).Value = 'Le Parlement se réunit en une session ordinaire'
$ws2.Cells.Item(10, 2).Value = 'De 9 mois'
$ws2.Cells.Item(10, 3).Value = 1
$ws2.Cells.Item(11, 2).Value = 'Permanente'
$ws2.Cells.Item(11, 3).Value = 0
$ws2.Cells.Item(12, 2).Value = 'De 8 mois'
$ws2.Cells.Item(12, 3).Value = 0
$ws2.Cells.Item(13, 2).Value = 'De 10 mois'
$ws2.Cells.Item(13, 3).Value = 0
$ws2.Cells.Item(14, 1).Value = 'Le Parlement exerce le pouvoir'
$ws2.Cells.Item(14, 2).Value = 'Législatif'
$ws2.Cells.Item(14, 3).Value = 1
$ws2.Cells.Item(15, 2).Value = 'Réglementaire'
$ws2.Cells.Item(15, 3).Value = 0
$ws2.Cells.Item(16, 2).Value = 'Administratif'
$ws2.Cells.Item(16, 3).Value = 0
$ws2.Cells.Item(17, 2).Value = 'Absolu'
$ws2.Cells.Item(17, 3).Value = 0
$ws2.Cells.Item(18, 1).Value = 'L''initiative des lois appartient, selon la Constitution,'
$ws2.Cells.Item(18, 2).Value = 'Concurremment au Premier ministre et aux membres du Parlement'
$ws2.Cells.Item(18, 3).Value = 1
$ws2.Cells.Item(19, 2).Value = 'Aux seuls membres du Parlement'
$ws2.Cells.Item(19, 3).Value = 0
$ws2.Cells.Item(20, 2).Value = 'Concurremment au Président de la République, au Premier ministre et aux membres du Parlement'
$ws2.Cells.Item(20, 3).Value = 0
$ws2.Cells.Item(21, 2).Value = 'Concurremment aux différents ministres et aux membres du Parlement'
$ws2.Cells.Item(21, 3).Value = 0
$ws2.Cells.Item(22, 1).Value = 'Depuis 1995, l''ordre du jour est arrêté par chaque Assemblée à l''initiative des groupes d''opposition ainsi qu''à celle des groupes minoritaires'
$ws2.Cells.Item(22, 2).Value = 'Un jour de séance par mois'
$ws2.Cells.Item(22, 3).Value = 1
$ws2.Cells.Item(23, 2).Value = 'Trois jours de séance par mois'
$ws2.Cells.Item(23, 3).Value = 0
$ws2.Cells.Item(24, 2).Value = 'Huit jours de séance par mois'
$ws2.Cells.Item(24, 3).Value = 0
$ws2.Cells.Item(25, 2).Value = 'Dix jours de séance par mois'
$ws2.Cells.Item(25, 3).Value = 0
$ws2.Cells.Item(26, 1).Value = 'Le mouvement d''un texte en discussion entre les députés et les sénateurs pour faire voter une loi s''appelle'
$ws2.Cells.Item(26, 2).Value = 'La navette parlementaire'
$ws2.Cells.Item(26, 3).Value = 1
$ws2.Cells.Item(27, 2).Value = 'Le cavalier législatif'
$ws2.Cells.Item(27, 3).Value = 0
$ws2.Cells.Item(28, 2).Value = 'La commission mixte paritaire parlementaire'
$ws2.Cells.Item(28, 3).Value = 0
$ws2.Cells.Item(29, 2).Value = 'Le va-et-vient parlementaire'
$ws2.Cells.Item(29, 3).Value = 0
$ws2.Cells.Item(30, 1).Value = 'Il y a 577 députés et 348 sénateurs qui seront prévus dêtre élus en 2011'
$ws2.Cells.Item(30, 3).Value = 1
$ws2.Cells.Item(31, 1).Value = 'Il y a 572 députés et 339 sénateurs qui seront prévus dêtre élus en 2011'
$ws2.Cells.Item(31, 3).Value = 0
$ws2.Cells.Item(32, 1).Value = 'Il y a 583 députés et 351 sénateurs qui seront prévus dêtre élus en 2011'
$ws2.Cells.Item(32, 3).Value = 0
$ws2.Cells.Item(33, 1).Value = 'Il y a 576 députés et 337 sénateurs qui seront prévus dêtre élus en 2011'
$ws2.Cells.Item(33, 3).Value = 0
$ws2.Cells.Item(34, 1).Value = 'Le mode d''élection des parlementaires est le scrutin uninominal à 2 tours pour les députés, le scrutin universel indirect pour les sénateurs'
$ws2.Cells.Item(34, 3).Value = 1
$ws2.Cells.Item(35, 1).Value = 'Le mode d''élection des parlementaires est le scrutin proportionnel plurinominal à 2 tours pour les députés, le scrutin universel direct pour les sénateurs'
$ws2.Cells.Item(35, 3).Value = 0
$ws2.Cells.Item(36, 1).Value = 'Le mode d''élection des parlementaires est le scrutin majoritaire plurinominal à 2 tours pour les députés, le scrutin universel direct pour les sénateurs'
$ws2.Cells.Item(36, 3).Value = 0
$ws2.Cells.Item(37, 1).Value = 'Le mode d''élection des parlementaires est le scrutin direct uninominal à 2 tours pour les députés, le scrutin universel indirect pour les sénateurs'
$ws2.Cells.Item(37, 3).Value = 0
$ws2.Cells.Item(38, 1).Value = 'Il faut avoir au moins 18 ans pour être député et au moins 24 ans pour être sénateur'
$ws2.Cells.Item(38, 3).Value = 1
$ws2.Cells.Item(39, 1).Value = 'Il faut avoir au moins 21 ans pour être député et au moins 28 ans pour être sénateur'
$ws2.Cells.Item(39, 3).Value = 0
$ws2.Cells.Item(40, 1).Value = 'Il faut avoir au moins 25 ans pour être député et au moins 32 ans pour être sénateur'
$ws2.Cells.Item(40, 3).Value = 0
$ws2.Cells.Item(41, 1).Value = 'Il faut avoir au moins 25 ans pour être député et au moins 35 ans pour être sénateur'
$ws2.Cells.Item(41, 3).Value = 0
$ws2.Cells.Item(42, 1).Value = 'Une inviolabilité pénale totale ne fait pas partie du statut du parlementaire'
$ws2.Cells.Item(42, 3).Value = 1
$ws2.Cells.Item(43, 1).Value = 'Une indemnité ne fait pas partie du statut du parlementaire'
$ws2.Cells.Item(43, 3).Value = 0
$ws2.Cells.Item(44, 1).Value = 'Une irresponsabilité parlementaire quant à ses opinions et son vote dans le cadre de son travail parlementaire ne fait pas partie du statut du parlementaire'
$ws2.Cells.Item(44, 3).Value = 0
$ws2.Cells.Item(45, 1).Value = 'Des incompatibilités de fonction ne fait pas partie du statut du parlementaire'
$ws2.Cells.Item(45, 3).Value = 0
$ws2.Cells.Item(46, 1).Value = 'La session ordinaire est ouverte et fermée par décret du Premier ministre'
$ws2.Cells.Item(46, 3).Value = 0
$ws2.Cells.Item(47, 1).Value = 'La session ordinaire est ouverte et fermée par décret du Président de la République'
$ws2.Cells.Item(47, 3).Value = 0
$ws2.Cells.Item(48, 1).Value = 'La session ordinaire est ouverte et fermée par arrêté conjoint des présidents des deux Assemblées'
$ws2.Cells.Item(48, 3).Value = 0
$ws2.Cells.Item(49, 1).Value = 'Il n''y a pas d''acte, c''est la Constitution qui fixe la période d''ouverture et de fermeture de la session ordinaire'
$ws2.Cells.Item(49, 3).Value = 1
$ws2.Cells.Item(50, 1).Value = 'Les sessions extraordinaires s''ouvrent et se ferment par décret du Président de la République'
$ws2.Cells.Item(50, 3).Value = 1
$ws2.Cells.Item(51, 1).Value = 'Les sessions extraordinaires s''ouvrent et se ferment par décret du Premier ministre'
$ws2.Cells.Item(51, 3).Value = 0
$ws2.Cells.Item(52, 1).Value = 'Les sessions extraordinaires s''ouvrent et se ferment par arrêté conjoint des présidents des deux Assemblées'
$ws2.Cells.Item(52, 3).Value = 0
$ws2.Cells.Item(53, 1).Value = 'Il n''y a pas d''acte, c''est la Constitution qui fixe la période d''ouverture et de fermeture des sessions extraordinaires'
$ws2.Cells.Item(53, 3).Value = 0
$ws2.Cells.Item(54, 1).Value = 'On peut trouver le compte rendu intégral des débats du Parlement selon la Constitution au Journal Officiel de la République française'
$ws2.Cells.Item(54, 3).Value = 1
$ws2.Cells.Item(55, 1).Value = 'On peut trouver le compte rendu intégral des débats du Parlement selon la Constitution sur le site internet de chacune des deux Assemblées (depuis 2003)'
$ws2.Cells.Item(55, 3).Value = 0
$ws2.Cells.Item(56, 1).Value = 'On peut trouver le compte rendu intégral des débats du Parlement selon la Constitution à la fois sur le site internet de chacune des deux Assemblées (depuis 2003) et au Journal Officiel de la République française'
$ws2.Cells.Item(56, 3).Value = 0
$ws2.Cells.Item(57, 1).Value = 'On peut trouver le compte rendu intégral des débats du Parlement selon la Constitution dans un numéro spécial de chaque ministère et qui paraît tous les 15 jours'
$ws2.Cells.Item(57, 3).Value = 0
$ws2.Cells.Item(58, 1).Value = 'La Constitution cite les domaines où la loi fixe les règles et ceux dont elle détermine les principes fondamentaux. La préservation de l''environnement relève des principes fondamentaux'
$ws2.Cells.Item(58, 3).Value = 1
$ws2.Cells.Item(59, 1).Value = 'La Constitution cite les domaines où la loi fixe les règles et ceux dont elle détermine les principes fondamentaux. La création de catégories d''établissements publics relève des principes fondamentaux'
$ws2.Cells.Item(59, 3).Value = 0
$ws2.Cells.Item(60, 1).Value = 'La Constitution cite les domaines où la loi fixe les règles et ceux dont elle détermine les principes fondamentaux. La nationalité, l''état et la capacité des personnes, les régimes matrimoniaux, les successions et libéralités relève des principes fondamentaux'
$ws2.Cells.Item(60, 3).Value = 0
$ws2.Cells.Item(61, 1).Value = 'La Constitution cite les domaines où la loi fixe les règles et ceux dont elle détermine les principes fondamentaux. L''assiette, le taux et les modalités de recouvrement des impositions de toutes natures ; le régime d''émission de la monnaie relève des principes fondamentaux'
$ws2.Cells.Item(61, 3).Value = 0
$ws2.Cells.Item(62, 1).Value = 'C''est Alain Poher qui assuma en 1969 et 1974 l''intérim de la Présidence de la République'
$ws2.Cells.Item(62, 3).Value = 1
$ws2.Cells.Item(63, 1).Value = 'C''est Gaston Monnerville qui assuma en 1969 et 1974 l''intérim de la Présidence de la République'
$ws2.Cells.Item(63, 3).Value = 0
$ws2.Cells.Item(64, 1).Value = 'C''est René Monory qui assuma en 1969 et 1974 l''intérim de la Présidence de la République'
$ws2.Cells.Item(64, 3).Value = 0
$ws2.Cells.Item(65, 1).Value = 'C''est Christian Poncelet qui assuma en 1969 et 1974 l''intérim de la Présidence de la République'
$ws2.Cells.Item(65, 3).Value = 0
$ws2.Cells.Item(66, 1).Value = '"Elles n''ont pas un caractère législatif dès leur publication". C''est une caractéristique des ordonnances que le Gouvernement peut prendre sur autorisation du Parlement dans un domaine qui est normalement du domaine de la loi '
$ws2.Cells.Item(66, 3).Value = 1
$ws2.Cells.Item(67, 1).Value = '"Elles sont prises en Conseil des Ministres après avis du Conseil d''État". C''est une caractéristique des ordonnances que le Gouvernement peut prendre sur autorisation du Parlement dans un domaine qui est normalement du domaine de la loi '
$ws2.Cells.Item(67, 3).Value = 0
$ws2.Cells.Item(68, 1).Value = '"Elles entrent en vigueur dès leur publication mais deviennent caduques si le projet de loi de ratification n''est pas déposé devant le Parlement". C''est une caractéristique des ordonnances que le Gouvernement peut prendre sur autorisation du Parlement dans un domaine qui est normalement du domaine de la loi .'
$ws2.Cells.Item(68, 3).Value = 0
$ws2.Cells.Item(69, 1).Value = '"Elles ne peuvent être ratifiées que de manière expresse". C''est une caractéristique des ordonnances que le Gouvernement peut prendre sur autorisation du Parlement dans un domaine qui est normalement du domaine de la loi '
$ws2.Cells.Item(69, 3).Value = 0
$ws2.Cells.Item(70, 1).Value = 'Relatif à la motion de censure, elle n''intervient que si le Président de l''Assemblée nationale donne son accord à son dépôt'
$ws2.Cells.Item(70, 3).Value = 1
$ws2.Cells.Item(71, 1).Value = 'Relatif à la motion de censure, elle est signée par un dixième au moins des membres de l''Assemblée nationale'
$ws2.Cells.Item(71, 3).Value = 0
$ws2.Cells.Item(72, 1).Value = 'Relatif à la motion de censure, le vote ne peut avoir lieu que quarante-huit heures après son dépôt'
$ws2.Cells.Item(72, 3).Value = 0
$ws2.Cells.Item(73, 1).Value = 'Relatif à la motion de censure, une majorité de votes favorables oblige le Premier ministre à remettre au Président de la République la démission du Gouvernement'
$ws2.Cells.Item(73, 3).Value = 0
$ws2.Cells.Item(74, 1).Value = 'Le parlement dispose de 70 jours pour voter la loi de finances (le budget) de l''année suivante'
$ws2.Cells.Item(74, 3).Value = 1
$ws2.Cells.Item(75, 1).Value = 'Le parlement dispose de 50 jours pour voter la loi de finances (le budget) de l''année suivante'
$ws2.Cells.Item(75, 3).Value = 0
$ws2.Cells.Item(76, 1).Value = 'Le parlement dispose de 60 jours pour voter la loi de finances (le budget) de l''année suivante'
$ws2.Cells.Item(76, 3).Value = 0
$ws2.Cells.Item(77, 1).Value = 'Le parlement dispose de 80 jours pour voter la loi de finances (le budget) de l''année suivante'
$ws2.Cells.Item(77, 3).Value = 0

# Column widths (auto-fit like Excel would do after pasting the data)
$ws2.Columns.Item(1).ColumnWidth = 255.83203125
$ws2.Columns.Item(2).ColumnWidth = 81.33203125
$ws2.Columns.Item(3).ColumnWidth = 9.33203125

# Leave the view positioned near the bottom of the newly-entered data,
# matching where the author ended up after typing in the new rows.
$ws2.Range("A78").Select()
